# Cache and checkpoint performance analysis (002).xlsx
# "Added command line for future" — add the 5,000,000-row benchmark column
# (and fill in a couple of previously-missing 6,000,000-row / NE-checkpoint
# points), rename Sheet1 to a descriptive name, and drop the legend from
# the second (stacked-area) chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Enter the new benchmark data -------------------------------------
# "No cache" / "Cache" / "Checkpoint" / "NE Checkpoint" timings for the
# 5,000,000-row run (column N) on the first table (rows 8-11).
$ws.Range("N8").Value = 102992
$ws.Range("N9").Value = 67581
$ws.Range("N10").Value = 50049
$ws.Range("N11").Value = 47707

# iMac 6,000,000-row run (column O) on the second table (rows 28-31).
$ws.Range("O28").Value = 88944
$ws.Range("O29").Value = 61033
$ws.Range("O30").Value = 45547
$ws.Range("O31").Value = 43707

# NE Checkpoint 5,000,000-row run (column N) on the third table (rows 33-36).
$ws.Range("N33").Value = 85851
$ws.Range("N34").Value = 56780
$ws.Range("N35").Value = 40071
$ws.Range("N36").Value = 33253

# --- 2. Rename the worksheet ----------------------------------------------
$oldName = $ws.Name
$newName = "Chapter 16 - Lab 100 Benchmark"
$ws.Name = $newName

# --- 3. Repoint the two charts' series formulas at the renamed sheet -----
# (rename doesn't cascade into chart SERIES() formulas automatically)
$chart1 = $ws.ChartObjects("Chart 1").Chart
$chart1.SeriesCollection(1).Formula = "=SERIES('$newName'!`$B`$48,'$newName'!`$H`$1:`$S`$1,'$newName'!`$H`$48:`$S`$48,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES('$newName'!`$B`$49,'$newName'!`$H`$1:`$S`$1,'$newName'!`$H`$49:`$S`$49,2)"
$chart1.SeriesCollection(3).Formula = "=SERIES('$newName'!`$B`$50,'$newName'!`$H`$1:`$S`$1,'$newName'!`$H`$50:`$S`$50,3)"
$chart1.SeriesCollection(4).Formula = "=SERIES('$newName'!`$B`$51,'$newName'!`$H`$1:`$S`$1,'$newName'!`$H`$51:`$S`$51,4)"

$chart2 = $ws.ChartObjects("Chart 2").Chart
$chart2.SeriesCollection(1).Formula = "=SERIES('$newName'!`$B`$43,'$newName'!`$C`$1:`$S`$1,'$newName'!`$C`$43:`$S`$43,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES('$newName'!`$B`$44,'$newName'!`$C`$1:`$S`$1,'$newName'!`$C`$44:`$S`$44,2)"
$chart2.SeriesCollection(3).Formula = "=SERIES('$newName'!`$B`$45,'$newName'!`$C`$1:`$S`$1,'$newName'!`$C`$45:`$S`$45,3)"
$chart2.SeriesCollection(4).Formula = "=SERIES('$newName'!`$B`$46,,'$newName'!`$C`$46:`$S`$46,4)"

# --- 4. Drop the legend from the second chart ------------------------------
$chart2.HasLegend = $false

# --- 5. Leave the selection on the last cell touched ----------------------
$ws.Range("M41").Select()
